$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their text formatting
# so values like "67.178.21" or "0.999" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '67.178.21'
$ws.Range("E2").Value = '  +0.46%  '

# Row 3
$ws.Range("D3").Value = '3.124.34'
$ws.Range("E3").Value = '  +0.76%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '580.11'
$ws.Range("E5").Value = '  -0.13%  '

# Row 6
$ws.Range("D6").Value = '174.59'
$ws.Range("E6").Value = '  +1.17%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("E8").Value = '  -0.23%  '

# Row 9
$ws.Range("B9").Value = 'Toncoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D9").Value = '6.45'
$ws.Range("E9").Value = '  +0.32%  '

# Row 10
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.155'
$ws.Range("E10").Value = '  +0.16%  '

# Row 11
$ws.Range("D11").Value = '0.481'
$ws.Range("E11").Value = '  -0.53%  '

# Row 12
$ws.Range("E12").Value = '  +0.15%  '

# Row 13
$ws.Range("E13").Value = '  -0.20%  '

# Row 15
$ws.Range("D15").Value = '3.639.31'
$ws.Range("E15").Value = '  +0.72%  '

# Row 16
$ws.Range("D16").Value = '67.145.22'
$ws.Range("E16").Value = '  +0.40%  '

# Row 17
$ws.Range("E17").Value = '  -0.55%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '16.72'
$ws.Range("E18").Value = '  +3.29%  '

# Row 19
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.120.96'
$ws.Range("E19").Value = '  +0.65%  '

# Row 20
$ws.Range("D20").Value = '490.13'
$ws.Range("E20").Value = '  +1.99%  '

# Row 21
$ws.Range("D21").Value = '7.93'
$ws.Range("E21").Value = '  +5.62%  '

# Row 22
$ws.Range("E22").Value = '  -0.84%  '

# Row 23
$ws.Range("D23").Value = '84.26'
$ws.Range("E23").Value = '  +0.38%  '

# Row 24
$ws.Range("D24").Value = '13.23'
$ws.Range("E24").Value = '  +0.34%  '

# Row 25
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").Value = '  -2.51%  '

# Row 26
$ws.Range("D26").Value = '10.44'
$ws.Range("E26").Value = '  +3.78%  '

# Row 27
$ws.Range("E27").Value = '  -0.04%  '

# Row 28
$ws.Range("D28").Value = '7.94'
$ws.Range("E28").Value = '  -0.26%  '

# Row 29
$ws.Range("E29").Value = '  -1.64%  '

# Row 30
$ws.Range("E30").Value = '  -0.43%  '

# Row 31
$ws.Range("D31").Value = '28.62'
$ws.Range("E31").Value = '  -0.61%  '

# Row 32
$ws.Range("E32").Value = '  -0.56%  '

# Row 33
$ws.Range("D33").Value = '0.0₃0952'
$ws.Range("E33").Value = '  -5.09%  '

# Row 34
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.08%  '

# Row 35
$ws.Range("D35").Value = '5.89'
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("E36").Value = '  -1.73%  '

# Row 37
$ws.Range("D37").Value = '47.34'
$ws.Range("E37").Value = '  -1.36%  '

# Row 38
$ws.Range("E38").Value = '  -2.67%  '

# Row 39
$ws.Range("D39").Value = '0.311'
$ws.Range("E39").Value = '  -1.51%  '

# Row 40
$ws.Range("D40").Value = '0.124'
$ws.Range("E40").Value = '  +1.82%  '

# Row 41
$ws.Range("D41").Value = '8.53'
$ws.Range("E41").Value = '  -1.49%  '

# Row 42
$ws.Range("D42").Value = '2.825.68'
$ws.Range("E42").Value = '  -0.44%  '

# Row 43
$ws.Range("D43").Value = '385.77'
$ws.Range("E43").Value = '  +0.10%  '

# Row 44
$ws.Range("D44").Value = '2.61'
$ws.Range("E44").Value = '  -7.04%  '

# Row 45
$ws.Range("E45").Value = '  -2.09%  '

# Row 46
$ws.Range("D46").Value = '135.53'
$ws.Range("E46").Value = '  +0.42%  '

# Row 47
$ws.Range("E47").Value = '  -0.02%  '

# Row 48
$ws.Range("E48").Value = '  -0.18%  '

# Row 49
$ws.Range("E49").Value = '  -1.17%  '

# Row 50
$ws.Range("E50").Value = '  -0.63%  '

# Row 51
$ws.Range("D51").Value = '6.76'
$ws.Range("E51").Value = '  -1.86%  '
